$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at K ("AccParticipation"), shifting the existing
#    Supported/STATUS/Key/Meaning columns one place to the right
#    (K->L, L->M, N->O, O->P).
$ws.Columns("K:K").Insert()
$ws.Range("K2").Value = "AccParticipation"

# 2. Append the four new account-level test-case rows (acc97-acc100).
$ws.Range("A99").Value = "acc97"
$ws.Range("D99").Value = '$'
$ws.Range("E99").Value = '$'
$ws.Range("F99").Value = '%'
$ws.Range("G99").Value = '%'
$ws.Range("H99").Value = '$'
$ws.Range("I99").Value = '$'
$ws.Range("J99").Value = '$'
$ws.Range("K99").Value = '%'

$ws.Range("A100").Value = "acc98"
$ws.Range("D100").Value = '$'
$ws.Range("E100").Value = '$'
$ws.Range("F100").Value = '%'
$ws.Range("G100").Value = '%'
$ws.Range("H100").Value = '$'
$ws.Range("I100").Value = '$'
$ws.Range("K100").Value = '%'

$ws.Range("A101").Value = "acc99"
$ws.Range("B101").Value = '$'
$ws.Range("D101").Value = '$'
$ws.Range("E101").Value = '$'
$ws.Range("F101").Value = '%'
$ws.Range("G101").Value = '%'
$ws.Range("H101").Value = '$'
$ws.Range("I101").Value = '$'
$ws.Range("J101").Value = '$'
$ws.Range("K101").Value = '%'

$ws.Range("A102").Value = "acc100"
$ws.Range("B102").Value = '$'
$ws.Range("D102").Value = '$'
$ws.Range("E102").Value = '$'
$ws.Range("F102").Value = '%'
$ws.Range("G102").Value = '%'
$ws.Range("H102").Value = '$'
$ws.Range("I102").Value = '$'
$ws.Range("K102").Value = '%'

# 3. Flip STATUS (now column M) from "In progress" to "complete" for every
#    test-case row, existing and newly-added alike.
for ($r = 3; $r -le 102; $r++) {
    $cell = $ws.Cells.Item($r, 13)
    $v = $cell.Value()
    if ($v -eq "In progress" -or $v -eq "") {
        $cell.Value = "complete"
    }
}

# 4. Restore the active cell selection to M4, as in the edited workbook.
$ws.Range("M4").Select()
